$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the Rafflesia Profits
# workbook: updates computed price/profit columns (H-N) across several rows
# in the ALC, CRP, CUL, GSM, LTW and WVR sheets.


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 80.833336
$ws.Range("I4").Value = 80.833336
$ws.Range("K4").Value = 80.833336
$ws.Range("M4").Value = 33.166664

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 825.8570999999999
$ws.Range("J33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("N33").Value = -2458

# Row 39: Riches' Brew
$ws.Range("H39").Value = 112.416664
$ws.Range("I39").Value = 61.555557
$ws.Range("J39").Value = 265
$ws.Range("K39").Value = 184.666671
$ws.Range("L39").Value = 795
$ws.Range("M39").Value = 111.333329
$ws.Range("N39").Value = -1387

# Row 48: The Sting of Conscience
$ws.Range("H48").Value = 7006.3335
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 7006.3335
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 21019.0005
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -21603.0005

# Row 56: Sleepless in Silvertear
$ws.Range("H56").Value = 7006.3335
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 7006.3335
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 21019.0005
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -22087.0005

# Row 58: A Matter of Vital Importance
$ws.Range("H58").Value = 1373.75
$ws.Range("I58").Value = 1600
$ws.Range("J58").Value = 1298.3334
$ws.Range("K58").Value = 4800
$ws.Range("L58").Value = 3895.0002
$ws.Range("M58").Value = -4650
$ws.Range("N58").Value = -4195.0002

# Row 64: Forged from the Void
$ws.Range("H64").Value = 5000.6665
$ws.Range("I64").Value = 5001
$ws.Range("K64").Value = 5001
$ws.Range("M64").Value = -4753

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 5000.6665
$ws.Range("I67").Value = 5001
$ws.Range("K67").Value = 5001
$ws.Range("M67").Value = -4143

# Row 95: Official Strategy Guide
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

# Row 97: Materia Worth
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 2250
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3582

# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 50
$ws.Range("I103").Value = 50
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 150
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 436
$ws.Range("N103").ClearContents()

# Row 106: Making Your Mark
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

# Row 113: Amaro Kart
$ws.Range("H113").Value = 1928.4286
$ws.Range("I113").Value = 1916.5
$ws.Range("K113").Value = 1916.5
$ws.Range("M113").Value = 1337.5

# Row 115: 5-bell Energy
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# Row 116: Growing Up
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

# Row 127: Liquid Competence
$ws.Range("H127").Value = 3532.3333
$ws.Range("I127").Value = 3438.8
$ws.Range("J127").Value = 4000
$ws.Range("K127").Value = 10316.4
$ws.Range("L127").Value = 12000
$ws.Range("M127").Value = -5356.400000000001
$ws.Range("N127").Value = -21920

# Row 131: Mindful Study
$ws.Range("H131").Value = 1303.8334
$ws.Range("I131").Value = 1303.8334
$ws.Range("K131").Value = 3911.5002
$ws.Range("M131").Value = 1128.4998


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 8000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -8406

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 5537.875
$ws.Range("I132").Value = 2078.4443
$ws.Range("K132").Value = 6235.3329
$ws.Range("M132").Value = -3705.3329

# Row 136: Turali Quality
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 24000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -29100


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 61: Red Letter Day
$ws.Range("H61").Value = 57
$ws.Range("I61").Value = 57
$ws.Range("K61").Value = 171
$ws.Range("M61").Value = 44

# Row 129: Comfort Food
$ws.Range("H129").Value = 4872.5
$ws.Range("I129").Value = 1998.3334
$ws.Range("J129").Value = 6597
$ws.Range("K129").Value = 5995.0002
$ws.Range("L129").Value = 19791
$ws.Range("M129").Value = -995.0002000000004
$ws.Range("N129").Value = -29791

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1336.6666
$ws.Range("I131").Value = 806
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 2418
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 2622
$ws.Range("N131").Value = -16080


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1581
$ws.Range("I102").Value = 1328.6666
$ws.Range("K102").Value = 1328.6666
$ws.Range("M102").Value = 293.3334


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 423.07693
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1090

# Row 27: Fire and Hide
$ws.Range("H27").Value = 423.07693
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -714

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 5500
$ws.Range("I93").Value = 5500
$ws.Range("K93").Value = 5500
$ws.Range("M93").Value = -4252

# Row 127: Loyal Turncoat
$ws.Range("H127").Value = 29455.092
$ws.Range("J127").Value = 29455.092
$ws.Range("L127").Value = 29455.092
$ws.Range("N127").Value = -39375.092


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5: Hire in the Blood
$ws.Range("H5").Value = 15003750
$ws.Range("J5").Value = 15003750
$ws.Range("L5").Value = 15003750
$ws.Range("N5").Value = -15003974

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2213
$ws.Range("I122").Value = 1426.5
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 4279.5
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -1829.5
$ws.Range("N122").Value = -13898.5
